$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Время начала заседания: [time]" -> "Время начала заседания: _____"
#    The underlined placeholder run "[time]" is replaced by literal
#    underscores and merged into the preceding (non-underlined) run.
# ---------------------------------------------------------------------------
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Text = "Время начала заседания: [time]"
$find.Replacement.Text = "Время начала заседания: _____"
$find.Forward = $true
$find.Wrap = 0
$find.Format = $false
$find.MatchCase = $false
$find.MatchWholeWord = $false
$find.MatchWildcards = $false
$find.MatchSoundsLike = $false
$find.MatchAllWordForms = $false
$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Collapse the triple space after the question-number parenthesis down to
#    a single space, e.g. ")   Краткое содержание вопроса" -> ") Краткое
#    содержание вопроса".
#
#    The very first such paragraph in the document keeps the numbering
#    ("1)   ") in its own run, separate from "Краткое содержание вопроса
#    ...", and that paragraph is NOT touched by this fix (matches the
#    original template / is left exactly as-is). Every later occurrence has
#    ")   Краткое содержание вопроса" inside a single run and gets the extra
#    spaces collapsed.
# ---------------------------------------------------------------------------
$marker = "Краткое содержание вопроса"
$seen = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $range = $p.Range
    $text = $range.Text

    if ($text.IndexOf($marker) -ge 0) {
        $seen = $seen + 1
        if ($seen -gt 1) {
            $pFind = $range.Find
            $pFind.ClearFormatting()
            $pFind.Replacement.ClearFormatting()
            $pFind.Text = ")   Краткое содержание вопроса"
            $pFind.Replacement.Text = ") Краткое содержание вопроса"
            $pFind.Forward = $true
            $pFind.Wrap = 0
            $pFind.Format = $false
            $pFind.MatchCase = $false
            $pFind.MatchWholeWord = $false
            $pFind.MatchWildcards = $false
            $pFind.MatchSoundsLike = $false
            $pFind.MatchAllWordForms = $false
            $pFind.Execute($pFind.Text, $false, $false, $false, $false, $false, $true, 1, $false, $pFind.Replacement.Text, 2) | Out-Null
        }
    }
}
